# Update crypto price/volume data per the GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value could be misread as a number (losing formatting
# like trailing zeros, or becoming numeric instead of text) are pinned to
# the Text format first, matching how the source data is stored.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "54.630.15"
$ws.Range("E2").Value = "  -7.30%  "
$ws.Range("D3").Value = "2.894.95"
$ws.Range("E3").Value = "  -10.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "478.16"
$ws.Range("E5").Value = "  -11.52%  "
$ws.Range("D6").Value = "127.33"
$ws.Range("E6").Value = "  -6.54%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "2.897.53"
$ws.Range("E8").Value = "  -10.18%  "
$ws.Range("D9").Value = "0.407"
$ws.Range("E9").Value = "  -11.41%  "
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  -11.85%  "
$ws.Range("E11").Value = "  -15.00%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -15.15%  "
$ws.Range("E13").Value = "  -3.85%  "
$ws.Range("D14").Value = "3.387.40"
$ws.Range("E14").Value = "  -10.38%  "
$ws.Range("D15").Value = "22.90"
$ws.Range("E15").Value = "  -11.97%  "
$ws.Range("D16").Value = "54.555.98"
$ws.Range("E16").Value = "  -7.49%  "
$ws.Range("D17").Value = "2.891.48"
$ws.Range("E17").Value = "  -10.40%  "
$ws.Range("E18").Value = "  -14.35%  "
$ws.Range("D19").Value = "5.28"
$ws.Range("E19").Value = "  -10.97%  "
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  -13.00%  "
$ws.Range("E21").Value = "  -13.38%  "
$ws.Range("D22").Value = "308.95"
$ws.Range("E22").Value = "  -14.93%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  -13.91%  "
$ws.Range("D25").Value = "59.75"
$ws.Range("E25").Value = "  -15.37%  "
$ws.Range("D26").Value = "0.996"
$ws.Range("D27").Value = "0.155"
$ws.Range("E27").Value = "  -9.64%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -15.15%  "
$ws.Range("D30").Value = "6.32"
$ws.Range("E30").Value = "  -11.40%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("D32").Value = "6.24"
$ws.Range("E32").Value = "  -12.20%  "
$ws.Range("D33").Value = "19.17"
$ws.Range("E33").Value = "  -12.52%  "
$ws.Range("E34").Value = "  -16.03%  "
$ws.Range("E35").Value = "  -13.81%  "
$ws.Range("D36").Value = "137.47"
$ws.Range("E36").Value = "  -14.99%  "
$ws.Range("E37").Value = "  -15.00%  "
$ws.Range("E38").Value = "  -15.48%  "
$ws.Range("D39").Value = "23.04"
$ws.Range("E39").Value = "  -12.70%  "
$ws.Range("E40").Value = "  -12.20%  "
$ws.Range("D41").Value = "2.920.95"
$ws.Range("E41").Value = "  -10.27%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "35.73"
$ws.Range("E43").Value = "  -13.18%  "
$ws.Range("D44").Value = "0.971"
$ws.Range("E44").Value = "  -12.59%  "
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").Value = "  -15.61%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.33"
$ws.Range("E46").Value = "  -11.76%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.44"
$ws.Range("E47").Value = "  -14.68%  "
$ws.Range("D48").Value = "2.059.55"
$ws.Range("E48").Value = "  -10.63%  "
$ws.Range("E49").Value = "  -15.49%  "
$ws.Range("D50").Value = "18.01"
$ws.Range("E50").Value = "  -13.60%  "
$ws.Range("D51").Value = "0.0214"
$ws.Range("E51").Value = "  -11.02%  "
